$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C (note) values for rows 2-14
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 9
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 8
$ws.Range("C10").Value = 3
$ws.Range("C11").Value = 8
$ws.Range("C12").Value = 5
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 5

# Column D (Commentaires) values - only specific rows get text.
# Order matches the original authoring order (shared-string table append order).
$ws.Range("D6").Value = "*****mauve sur noir n'est pas lisible"
$ws.Range("D7").Value = "* ID's sont unique et ne peuvent pas être utiliser plus qu'une fois dans le document"
$ws.Range("D10").Value = "voir commentaire de navigation"
$ws.Range("D11").Value = "confusion sur quel document à coriger, vue que deux fichier de projets on été remis"
$ws.Range("D12").Value = "description n'as pas été modifier de stock"
$ws.Range("D13").Value = "aucun commentaires"
$ws.Range("D9").Value = "erreurs dans le CSS ** balises <style> ne vont pas dans un fichier .css"

# Apply the wrap-text style (same as column A) to D2:D15
$ws.Range("D2:D15").WrapText = $true

# Update the active selection to D10 as in the saved file
$ws.Range("D10").Select()

$wb.Save()
